# Update functions and Data Model (#50)
#
# Adds a new "Authorship Resource" column (14th column, N) to Table1,
# labels the header and fills every data row (2-52) with the
# author-credit string, then nudges the view/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Add the new table column ------------------------------------------------
$col = $lo.ListColumns.Add()
$col.Range.Item(1).Value = "Authorship Resource"

# Fill the whole data body in one shot, then make sure the applied format
# matches the rest of column N (font / vertical alignment) so it resolves
# to the same cell style already used by the existing N2:N16 cells.
$col.DataBodyRange.Value = "Daniela Subotic, Noémi Villars-Amberg"
$col.DataBodyRange.Font.Name = "Arial"
$col.DataBodyRange.Font.Size = 14
$col.DataBodyRange.VerticalAlignment = -4160

# --- View / selection ---------------------------------------------------------
# Re-establish the header freeze (row 1) and leave the selection on the
# newly authored column, matching the saved selection in the workbook.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("N2:N52").Select()
